$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "27.596.29"
$c.Style = "Normal"
$ws.Range("E2").Value = "  -1.51%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.667.99"
$c.Style = "Normal"
$ws.Range("E3").Value = "  -3.03%  "
$ws.Range("E4").Value = "  -0.16%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "215.09"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -1.65%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.511"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -2.20%  "
$ws.Range("E7").Value = "  -0.13%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "23.69"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -1.82%  "
$ws.Range("E9").Value = "  -0.71%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.0622"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -1.48%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0881"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -1.89%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "1.904.51"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -3.03%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "1.678.80"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -2.40%  "
$ws.Range("E14").Value = "  -2.99%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.561"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +0.25%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "66.24"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -1.69%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "27.588.27"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -1.41%  "
$ws.Range("B18").Value = "BitcoinCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "243.22"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +0.49%  "
$ws.Range("E19").Value = "  -3.28%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "7.62"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -3.68%  "
$ws.Range("E22").Value = "  -2.83%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "9.29"
$c.Style = "Normal"
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "2.03"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -4.53%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "146.88"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -1.21%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "7.20"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -3.78%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "16.44"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -1.33%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("E29").Value = "  -2.18%  "
$ws.Range("E30").Value = "  +2.96%  "
$ws.Range("E31").Value = "  -1.36%  "
$ws.Range("E32").Value = "  -2.38%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "1.465.99"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -1.35%  "
$ws.Range("E34").Value = "  -4.55%  "
$ws.Range("E35").Value = "  -4.88%  "
$ws.Range("E36").Value = "  -1.70%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.928"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -2.51%  "
$ws.Range("E38").Value = "  -0.95%  "
$ws.Range("E39").Value = "  -5.22%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "69.51"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -1.38%  "
$ws.Range("E41").Value = "  -5.36%  "
$ws.Range("E43").Value = "  -6.65%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "2.23"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -2.93%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "1.811.56"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -3.06%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.788"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -0.95%  "
$ws.Range("E47").Value = "  -1.87%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "89.36"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -1.59%  "
$ws.Range("E49").Value = "  -3.53%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "7.90"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -4.20%  "
